$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.921.58'
$ws.Range('E2').Value = '  +0.59%  '
$ws.Range('D3').Value = '1.674.00'
$ws.Range('E3').Value = '  +2.29%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '216.11'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.91%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.533'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +5.81%  '
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('E8').Value = '  +3.11%  '
$ws.Range('E9').Value = '  +1.68%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.37'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +4.58%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0892'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.97%  '
$ws.Range('D12').Value = '1.909.17'
$ws.Range('E12').Value = '  +2.13%  '
$ws.Range('D13').Value = '1.665.70'
$ws.Range('E13').Value = '  +1.60%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.10'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.03%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.525'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.09%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.75'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.92%  '
$ws.Range('D17').Value = '26.939.85'
$ws.Range('E17').Value = '  +0.58%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '233.06'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -3.57%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.90'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.84%  '
$ws.Range('E20').Value = '  +1.46%  '
$ws.Range('E21').Value = '  +0.22%  '
$ws.Range('E22').Value = '  +2.73%  '
$ws.Range('B23').Value = 'Toncoin'
$ws.Range('C23').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.22'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.54%  '
$ws.Range('B24').Value = 'Avalanche'
$ws.Range('C24').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.22'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.06%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '145.71'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.11%  '
$ws.Range('E26').Value = '  +1.34%  '
$ws.Range('E27').Value = '  +2.26%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.98'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.19%  '
$ws.Range('E29').Value = '  +0.10%  '
$ws.Range('E30').Value = '  +0.53%  '
$ws.Range('E31').Value = '  +1.31%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.34'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.10%  '
$ws.Range('D33').Value = '1.465.17'
$ws.Range('E33').Value = '  -2.98%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.18'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +4.94%  '
$ws.Range('E35').Value = '  +5.07%  '
$ws.Range('E36').Value = '  -0.07%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.908'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +5.66%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.571'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.54%  '
$ws.Range('E39').Value = '  +0.84%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.07'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.43%  '
$ws.Range('E41').Value = '  +0.22%  '
$ws.Range('E42').Value = '  +4.54%  '
$ws.Range('B43').Value = 'WEMIXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.974'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +6.91%  '
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '65.89'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.90%  '
$ws.Range('D45').Value = '1.820.18'
$ws.Range('E45').Value = '  +2.23%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.780'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.53%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '90.67'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.62%  '
$ws.Range('E48').Value = '  +1.06%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').Value = '0.0₆0102'
$ws.Range('E49').Value = '  +12.44%  '
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.100'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.83%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0508'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.39%  '
